$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excluded structures")

# Fill in row 17 with test data.
$ws.Range("A17").Value = "for testing"
$ws.Range("B17").Value = "for testing"
$ws.Range("C17").Value = "for testing"
$ws.Range("D17").Value = "for testing"
$ws.Range("E17").Value = "for testing"
$ws.Range("F17").Value = "for testing"
$ws.Range("G17").Value = "for testing"
$ws.Range("H17").Value = "No structure"
$ws.Range("I17").Value = "Imagery review"
$ws.Range("J17").Value = "for testing"
$ws.Range("K17").Value = "for testing"
$ws.Range("J17:K17").WrapText = $true

$ws.Range("K17").Select()
